# Add a new "raw meat" food item as row 24 on the "items" sheet.
# Write order matters: new shared-string entries are appended to the
# sharedStrings table in the order the cell values are first assigned, and
# the target workbook expects them in the order
#   itd_rawmeat, 200, 50, it_fo_rawmeat
# which corresponds to writing column B, then D, then A before the rest.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B24").Value = "itd_rawmeat"
$ws.Range("D24").Value = "200, 50"
$ws.Range("A24").Value = "it_fo_rawmeat"
$ws.Range("C24").Value = 8
$ws.Range("E24").Value = 2
$ws.Range("G24").Value = 2
$ws.Range("J24").Value = 5
$ws.Range("T24").Value = 10

# Move/select the cell the author's cursor ended up on after the edit.
$ws.Range("C25").Select() | Out-Null
